# Generate Report for Handback
# The f03ca9ea-... file has now been handed back (in sync with en-US),
# so its status moves from "Ready for handoff" to
# "Handed back: in sync with en-US", and the "Latest Handback DateTime"
# is refreshed to reflect when the handback report was generated.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 corresponds to f03ca9ea-...md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet: row 3 corresponds to f03ca9ea-...md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = $newStatus
$wsZhCn.Range("G2").Value = "2016-03-09 16:14:06"
$wsZhCn.Range("G3").Value = "2016-03-09 16:14:06"

# --- de-de sheet: row 3 corresponds to f03ca9ea-...md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = $newStatus
$wsDeDe.Range("G2").Value = "2016-03-09 16:14:22"
$wsDeDe.Range("G3").Value = "2016-03-09 16:14:22"
